# Apply the dated-worksheet update: refresh the date heading and all
# division problems in the practice table to the new "aa3dc9e" values.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Heading date.
Replace-Text "2023-11-22 Wednesday" "2023-11-23 Thursday"

# Division problems that only occur once in the document can be replaced
# via a simple Find/Replace.
Replace-Text "84÷4=" "57÷9="
Replace-Text "22÷3=" "67÷8="
Replace-Text "62÷4=" "42÷6="
Replace-Text "12÷2=" "76÷7="
Replace-Text "97÷4=" "21÷8="
Replace-Text "73÷6=" "58÷2="
Replace-Text "83÷3=" "58÷4="
Replace-Text "56÷6=" "32÷6="
Replace-Text "69÷3=" "86÷4="
Replace-Text "39÷2=" "47÷9="
Replace-Text "52÷5=" "17÷6="
Replace-Text "40÷4=" "13÷3="
Replace-Text "43÷8=" "55÷6="
Replace-Text "97÷9=" "18÷5="
Replace-Text "79÷4=" "11÷4="
Replace-Text "44÷7=" "39÷3="
Replace-Text "12÷9=" "21÷6="
Replace-Text "78÷7=" "90÷9="
Replace-Text "86÷8=" "11÷9="
Replace-Text "91÷6=" "79÷5="
Replace-Text "75÷3=" "81÷7="
Replace-Text "46÷5=" "79÷9="
Replace-Text "81÷8=" "84÷9="

# "47÷7=" occurs twice in the table (row 1, col 5 and row 13, col 5) and
# each instance maps to a different new value, so address those two
# cells directly by position instead of relying on Find/Replace.
$table = $d.Tables.Item(1)
$table.Cell(1, 5).Range.Text = "10÷8="
$table.Cell(13, 5).Range.Text = "28÷4="
